# edit.ps1
# Applies the cryptos.xlsx data refresh described in the commit:
# "Updated cryptos list on Sat May 13 03:56:54 UTC 2023 with GitHub Actions"
#
# For each data row (2-51), update the Price (D) and Volume(1h) (E) text
# values with the freshly scraped figures. For rows 29/30 and 32/33 the
# underlying coins were also reordered (swapped), so Coin (B) and Link (C)
# are updated there too.
#
# Price values are stored as plain text in the source sheet (e.g. thousands
# separated by '.', like "26.916.11"), so assignments are done with a
# leading apostrophe to keep Excel from reinterpreting numeric-looking
# strings (e.g. "1.003") as numbers; the cell style is reset back to
# "Normal" immediately after so no stray number-format/quote-prefix style
# is left attached to the cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Row = 2; D='26.916.11'; E='  +0.38%  ' }
    @{ Row = 3; D='1.813.24'; E='  +1.88%  ' }
    @{ Row = 4; D='1.003'; E='  -0.73%  ' }
    @{ Row = 5; D='311.89'; E='  +1.46%  ' }
    @{ Row = 6; E='  -0.46%  ' }
    @{ Row = 7; D='0.4284'; E='  +1.72%  ' }
    @{ Row = 8; E='  +2.47%  ' }
    @{ Row = 9; D='0.07240'; E='  +1.22%  ' }
    @{ Row = 10; D='0.8629'; E='  +3.46%  ' }
    @{ Row = 11; D='21.24'; E='  +5.23%  ' }
    @{ Row = 12; D='2.008.69'; E='  +10.21%  ' }
    @{ Row = 13; D='6.627'; E='  +4.89%  ' }
    @{ Row = 14; D='5.389'; E='  +3.00%  ' }
    @{ Row = 15; D='0.06918'; E='  +1.66%  ' }
    @{ Row = 16; D='80.63'; E='  +2.07%  ' }
    @{ Row = 17; E='  -0.90%  ' }
    @{ Row = 18; D='0.000008926'; E='  +3.05%  ' }
    @{ Row = 19; E='  -0.41%  ' }
    @{ Row = 20; E='  +2.28%  ' }
    @{ Row = 21; D='26.943.72'; E='  -0.06%  ' }
    @{ Row = 22; D='5.177'; E='  +3.27%  ' }
    @{ Row = 23; D='11.02'; E='  +0.15%  ' }
    @{ Row = 24; D='2.232.78'; E='  +9.77%  ' }
    @{ Row = 25; D='153.59'; E='  +0.49%  ' }
    @{ Row = 26; D='1.885'; E='  -1.55%  ' }
    @{ Row = 27; E='  +1.02%  ' }
    @{ Row = 28; D='5.234'; E='  +4.50%  ' }
    @{ Row = 29; B='LidoDAOToken'; C='https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'; D='1.883'; E='  +16.19%  ' }
    @{ Row = 30; B='BitcoinCash'; C='https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'; D='114.81'; E='  +0.48%  ' }
    @{ Row = 31; D='0.08954'; E='  +0.30%  ' }
    @{ Row = 32; B='ImmutableX'; C='https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'; D='0.7434'; E='  +4.80%  ' }
    @{ Row = 33; B='ARBITRUM'; C='https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'; D='1.165'; E='  +7.37%  ' }
    @{ Row = 34; D='4.433'; E='  +3.30%  ' }
    @{ Row = 35; D='2.798'; E='  -1.57%  ' }
    @{ Row = 36; D='1.007'; E='  -0.05%  ' }
    @{ Row = 37; D='1.115'; E='  +3.58%  ' }
    @{ Row = 38; D='0.05213'; E='  +2.85%  ' }
    @{ Row = 39; D='0.01921'; E='  +2.01%  ' }
    @{ Row = 40; D='0.5084'; E='  +3.77%  ' }
    @{ Row = 41; D='0.1642'; E='  +2.44%  ' }
    @{ Row = 42; D='2.723'; E='  +8.01%  ' }
    @{ Row = 43; D='6.442'; E='  +8.06%  ' }
    @{ Row = 44; D='8.271'; E='  +4.97%  ' }
    @{ Row = 45; D='106.83'; E='  +2.70%  ' }
    @{ Row = 46; D='10.43'; E='  +3.40%  ' }
    @{ Row = 47; D='1.004'; E='  -0.43%  ' }
    @{ Row = 48; D='1.654'; E='  +5.58%  ' }
    @{ Row = 49; D='0.4578'; E='  +3.00%  ' }
    @{ Row = 50; D='0.06279'; E='  +0.55%  ' }
    @{ Row = 51; D='1.804'; E='  +6.53%  ' }
)

foreach ($u in $updates) {
    $r = $u.Row
    if ($u.ContainsKey('B')) {
        $ws.Cells.Item($r, 2).Value = $u.B
    }
    if ($u.ContainsKey('C')) {
        $ws.Cells.Item($r, 3).Value = $u.C
    }
    if ($u.ContainsKey('D')) {
        $cell = $ws.Cells.Item($r, 4)
        $cell.Value = "'" + $u.D
        $cell.Style = "Normal"
    }
    if ($u.ContainsKey('E')) {
        $ws.Cells.Item($r, 5).Value = $u.E
    }
}
